# Auto update Excel log
# Appends new sensor-log rows (captured 2026-01-28, ~14:43-14:53) to the
# PIR, Humidity, Temperature and mmWave sheets of the SeniorConnect master
# log workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write one data row (Date, Timestamp, Hour, Location, Value,
# Status) into row $r of worksheet $ws, columns A:F. Column A (the date
# string) - and, when requested, column E - are quote-prefixed so Excel's
# automatic type-sniffing on Range.Value keeps them as literal text
# instead of silently coercing "2026-01-28" to a date serial or "88.9%"
# to a percentage number.
# ---------------------------------------------------------------------
function Add-LogRow {
    param($ws, $r, $date, $timestamp, $hour, $location, $value, $status, $forceTextValue)

    $ws.Cells.Item($r, 1).Value = "'" + $date
    $ws.Cells.Item($r, 2).Value = $timestamp
    $ws.Cells.Item($r, 3).Value = $hour
    $ws.Cells.Item($r, 4).Value = $location
    if ($forceTextValue) {
        $ws.Cells.Item($r, 5).Value = "'" + $value
    } else {
        $ws.Cells.Item($r, 5).Value = $value
    }
    $ws.Cells.Item($r, 6).Value = $status
}

# ---------------------------------------------------------------------
# PIR sheet: append rows 53-70
# ---------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-01-28","14:43:50","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:43:53","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:43:58","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:44:03","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:44:08","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:44:13","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:52:22","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:52:27","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:52:32","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:52:37","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:52:42","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:52:47","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:52:52","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:52:57","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:53:02","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:53:07","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:53:12","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:53:17","14:00","Bathroom","No Motion","Inactive")
)
$r = 53
foreach ($row in $pirRows) {
    Add-LogRow $pir $r $row[0] $row[1] $row[2] $row[3] $row[4] $row[5] $false
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Humidity sheet: append rows 49-66 (Value column is a "%" reading, also
# quote-prefixed so it is preserved as text, not parsed as a percentage)
# ---------------------------------------------------------------------
$humidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("2026-01-28","14:43:56","14:00","Bathroom","88.9%","Active"),
    @("2026-01-28","14:44:00","14:00","Bathroom","88.9%","Active"),
    @("2026-01-28","14:44:04","14:00","Bathroom","88.0%","Active"),
    @("2026-01-28","14:44:08","14:00","Bathroom","88.8%","Active"),
    @("2026-01-28","14:44:12","14:00","Bathroom","87.8%","Active"),
    @("2026-01-28","14:52:19","14:00","Bathroom","88.3%","Active"),
    @("2026-01-28","14:52:23","14:00","Bathroom","87.4%","Active"),
    @("2026-01-28","14:52:27","14:00","Bathroom","88.3%","Active"),
    @("2026-01-28","14:52:31","14:00","Bathroom","88.3%","Active"),
    @("2026-01-28","14:52:35","14:00","Bathroom","87.4%","Active"),
    @("2026-01-28","14:52:39","14:00","Bathroom","88.4%","Active"),
    @("2026-01-28","14:52:47","14:00","Bathroom","87.4%","Active"),
    @("2026-01-28","14:52:52","14:00","Bathroom","88.3%","Active"),
    @("2026-01-28","14:52:55","14:00","Bathroom","87.4%","Active"),
    @("2026-01-28","14:53:00","14:00","Bathroom","88.3%","Active"),
    @("2026-01-28","14:53:08","14:00","Bathroom","88.3%","Active"),
    @("2026-01-28","14:53:12","14:00","Bathroom","88.3%","Active"),
    @("2026-01-28","14:53:16","14:00","Bathroom","87.4%","Active")
)
$r = 49
foreach ($row in $humidityRows) {
    Add-LogRow $humidity $r $row[0] $row[1] $row[2] $row[3] $row[4] $row[5] $true
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Temperature sheet: append rows 49-66
# ---------------------------------------------------------------------
$temperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("2026-01-28","14:43:57","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:44:01","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:44:05","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:44:09","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:44:13","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:52:19","14:00","Bathroom","22.7C","Active"),
    @("2026-01-28","14:52:24","14:00","Bathroom","22.8C","Active"),
    @("2026-01-28","14:52:28","14:00","Bathroom","22.8C","Active"),
    @("2026-01-28","14:52:32","14:00","Bathroom","22.8C","Active"),
    @("2026-01-28","14:52:36","14:00","Bathroom","22.8C","Active"),
    @("2026-01-28","14:52:40","14:00","Bathroom","22.8C","Active"),
    @("2026-01-28","14:52:48","14:00","Bathroom","22.8C","Active"),
    @("2026-01-28","14:52:52","14:00","Bathroom","22.8C","Active"),
    @("2026-01-28","14:52:56","14:00","Bathroom","22.8C","Active"),
    @("2026-01-28","14:53:00","14:00","Bathroom","22.8C","Active"),
    @("2026-01-28","14:53:08","14:00","Bathroom","22.8C","Active"),
    @("2026-01-28","14:53:12","14:00","Bathroom","22.8C","Active"),
    @("2026-01-28","14:53:16","14:00","Bathroom","22.8C","Active")
)
$r = 49
foreach ($row in $temperatureRows) {
    Add-LogRow $temperature $r $row[0] $row[1] $row[2] $row[3] $row[4] $row[5] $false
    $r = $r + 1
}

# ---------------------------------------------------------------------
# mmWave sheet: append row 3
# ---------------------------------------------------------------------
$mmwave = $wb.Worksheets.Item("mmWave")
Add-LogRow $mmwave 3 "2026-01-28" "14:53:11" "14:00" "Living Room" "No Presence" "Inactive" $false
